# Recitation_12.pptx - slide 18 "Step 2" caption textbox:
#  - reposition/resize the textbox (it becomes a single wide line instead of
#    a narrow multi-line box)
#  - split the trailing sentence so "beginning of the FREE LIST (LIFO)" is
#    called out in bold red

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(18)
$sh = $s.Shapes.Item(15)

# --- reposition / resize -------------------------------------------------
# Target EMU values (from the authored OOXML): off x=5561814 y=3901556,
# ext cx=6240546 cy=369332. Shape.Left/Top/Width/Height are Single (f32)
# COM properties, so we hand them point values that round-trip to the exact
# EMU figures above rather than the naively-rounded EMU/12700 value.
$sh.Left   = 437.9381103515625
$sh.Width  = 491.381591796875
$sh.Height = 29.081260681152344
# Top (y) is unchanged from the original, so it is left untouched.

# --- split the second run -------------------------------------------------
$tr   = $sh.TextFrame.TextRange
$run2 = $tr.Runs(2)
$run2.Text = " Insert Block 4.2 to the "
$null = $run2.InsertAfter("beginning of the FREE LIST (LIFO)")

# Re-fetch the freshly split-off third run and format it as bold red.
$run3 = $sh.TextFrame.TextRange.Runs(3)
$run3.Font.Bold = $true
$run3.Font.Color.RGB = 255
